$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.968.76'
$ws.Range('E2').Value = '  +5.02%  '

$ws.Range('D3').Value = '2.243.56'
$ws.Range('E3').Value = '  +2.25%  '

$ws.Range('E4').Value = '  -0.21%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.22'
$ws.Range('E5').Value = '  +2.30%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.629'
$ws.Range('E6').Value = '  +0.15%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '61.77'
$ws.Range('E7').Value = '  -2.34%  '

$ws.Range('E8').Value = '  -0.07%  '

$ws.Range('E9').Value = '  +3.07%  '

$ws.Range('E10').Value = '  +1.31%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0922'
$ws.Range('E11').Value = '  +7.47%  '

$ws.Range('E12').Value = '  +0.74%  '

$ws.Range('D13').Value = '2.576.03'
$ws.Range('E13').Value = '  +2.24%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.81'
$ws.Range('E14').Value = '  +0.32%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '22.59'
$ws.Range('E15').Value = '  +2.77%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.811'
$ws.Range('E16').Value = '  -0.51%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.64'
$ws.Range('E17').Value = '  +1.52%  '

$ws.Range('D18').Value = '2.244.20'
$ws.Range('E18').Value = '  +2.36%  '

$ws.Range('D19').Value = '41.913.79'
$ws.Range('E19').Value = '  +4.91%  '

$ws.Range('D20').Value = '0.0₃0923'
$ws.Range('E20').Value = '  +2.13%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.80'
$ws.Range('E21').Value = '  +0.80%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.10'
$ws.Range('E22').Value = '  +0.59%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '253.70'
$ws.Range('E23').Value = '  +9.03%  '

$ws.Range('E24').Value = '  -0.08%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.40'
$ws.Range('E25').Value = '  +2.14%  '

$ws.Range('E26').Value = '  -0.34%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.80'
$ws.Range('E27').Value = '  +1.60%  '

$ws.Range('E28').Value = '  +3.88%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '169.82'
$ws.Range('E29').Value = '  -1.27%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.19'
$ws.Range('E30').Value = '  +0.42%  '

$ws.Range('E31').Value = '  -1.19%  '

$ws.Range('E32').Value = '  -0.37%  '

$ws.Range('E33').Value = '  -0.03%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.11'
$ws.Range('E34').Value = '  +8.75%  '

$ws.Range('E35').Value = '  +3.23%  '

$ws.Range('E36').Value = '  +2.68%  '

$ws.Range('E37').Value = '  -3.93%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.73'
$ws.Range('E38').Value = '  -3.43%  '

$ws.Range('E39').Value = '  -2.11%  '

$ws.Range('E40').Value = '  +30.97%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  -0.18%  '

$ws.Range('E42').Value = '  +6.30%  '

$ws.Range('E43').Value = '  +4.48%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.60'
$ws.Range('E44').Value = '  -8.23%  '

$ws.Range('E45').Value = '  +0.77%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '100.12'
$ws.Range('E46').Value = '  -2.69%  '

$ws.Range('E47').Value = '  +3.31%  '

$ws.Range('D48').Value = '1.485.46'
$ws.Range('E48').Value = '  -1.95%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '16.77'
$ws.Range('E49').Value = '  -3.88%  '

$ws.Range('E50').Value = '  +0.04%  '

$ws.Range('B51').Value = 'ARBITRUM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.09'
$ws.Range('E51').Value = '  -0.76%  '
